$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row label changes (shared-string reorder: Gipuzkoa/Guipuzcoa moved up) ---
$ws.Range("A22").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("A23").Value = "Sevilla"
$ws.Range("A24").Value = "Asturias"

# --- Updated statistics (Casos totales, Casos activos, Recuperados, Muertes) ---

# Row 7 (Bizkaia/Vizcaya)
$ws.Range("B7").Value = 5798
$ws.Range("C7").Value = 4867
$ws.Range("D7").Value = 5101
$ws.Range("E7").Value = 414

# Row 16 (Araba/Alava)
$ws.Range("B16").Value = 3019
$ws.Range("C16").Value = 4867
$ws.Range("D16").Value = 5101
$ws.Range("E16").Value = 260

# Row 22 (now Gipuzkoa/Guipuzcoa)
$ws.Range("B22").Value = 1955
$ws.Range("C22").Value = 4867
$ws.Range("D22").Value = 5101
$ws.Range("E22").Value = 130

# Row 23 (now Sevilla)
$ws.Range("B23").Value = 1947
$ws.Range("C23").Value = 221
$ws.Range("D23").Value = 1567
$ws.Range("E23").Value = 159

# Row 24 (now Asturias)
$ws.Range("B24").Value = 1892
$ws.Range("C24").Value = 430
$ws.Range("D24").Value = 1322
$ws.Range("E24").Value = 140

# Row 28 (Caceres)
$ws.Range("B28").Value = 1721
$ws.Range("C28").Value = 229
$ws.Range("D28").Value = 1235
$ws.Range("E28").Value = 257

# Row 43 (Avila)
$ws.Range("B43").Value = 859
$ws.Range("C43").Value = 252
$ws.Range("D43").Value = 552
$ws.Range("E43").Value = 55

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 13:22"
